$d = $word.ActiveDocument

# --- Paragraph 1 (title): "Kur-Beschreibung ... Routenplanung..." -> "Kurzbeschreibung ... Routenplanung..." (single line) ---
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("Kur-Beschreibung der Arbeitsweise im Projekt Routenplanung mit natürlicher Sprache", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Kurzbeschreibung der Arbeitsweise im Projekt Routenplanung mit natürlicher Sprache", 2)

# --- Paragraph 2: merge the three runs quoting the project name into flowing text (same content) ---
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("Im Folgenden soll die Arbeitsweise im Projekt „Routenplanung mit natürlicher Sprache“ kurz erläutert werden. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Im Folgenden soll die Arbeitsweise im Projekt „Routenplanung mit natürlicher Sprache“ kurz erläutert werden. ", 2)

# --- Paragraph 3: "das" -> "dass" ---
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("verständigt, das zunächst", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "verständigt, dass zunächst", 2)

# --- Paragraph 6 (Nachdem Anfang April...): several fixes + new content ---
$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute("konzentrierten sich", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "konzentrierte sich", 2)

$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute("Literatur-Recherche,", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Literaturrecherche,", 2)

$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute("Kosten-Zeit-Aufwandschätzung sowie", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Kosten-Zeit-Aufwandschätzung, die Erstellung des visuellen Prototyps sowie", 2)

$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute("gewährleisten. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "gewährleisten. Hierbei wurde die Customer Journey als Basis verwendet, um den Nutzen für den Kunden darzustellen.", 2)

# --- Paragraph 7 (Ende Juni...): "Schreib-prozess" -> "Schreibprozess" ---
$p7 = $d.Paragraphs(7).Range
$p7.Find.Execute("Schreib-prozess", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Schreibprozess", 2)

# --- Paragraph 8 (Elementar...): "v.a." -> "v. a." and "Github" -> "GitHub" ---
$p8 = $d.Paragraphs(8).Range
$p8.Find.Execute("dabei v.a. der", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "dabei v. a. der", 2)

$p8 = $d.Paragraphs(8).Range
$p8.Find.Execute("mittels Github.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "mittels GitHub.", 2)
